# Auto-generated edit script: updates numeric leve-profit cells per the commit diff.
# All target cells are plain numeric values (no formulas in the source workbook).
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5992.5
$ws.Range("I40").Value = 1985
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 1985
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -1810
$ws.Range("N40").Value = -10350
# Row 51
$ws.Range("H51").Value = 7101.3335
$ws.Range("J51").Value = 7101.3335
$ws.Range("L51").Value = 7101.3335
$ws.Range("N51").Value = -8069.3335
# Row 64
$ws.Range("H64").Value = 5532.3335
$ws.Range("J64").Value = 7333
$ws.Range("L64").Value = 7333
$ws.Range("N64").Value = -7829
# Row 67
$ws.Range("H67").Value = 5532.3335
$ws.Range("J67").Value = 7333
$ws.Range("L67").Value = 7333
$ws.Range("N67").Value = -9049
# Row 76
$ws.Range("H76").Value = 5538.04
$ws.Range("I76").Value = 5389.5264
$ws.Range("J76").Value = 6008.3335
$ws.Range("K76").Value = 5389.5264
$ws.Range("L76").Value = 6008.3335
$ws.Range("M76").Value = -5074.5264
$ws.Range("N76").Value = -6638.3335
# Row 79
$ws.Range("H79").Value = 5538.04
$ws.Range("I79").Value = 5389.5264
$ws.Range("J79").Value = 6008.3335
$ws.Range("K79").Value = 5389.5264
$ws.Range("L79").Value = 6008.3335
$ws.Range("M79").Value = -4297.5264
$ws.Range("N79").Value = -8192.333500000001
# Row 100
$ws.Range("H100").Value = 2087.2144
$ws.Range("I100").Value = 1851.75
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 1851.75
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -1310.75
$ws.Range("N100").Value = -4582
# Row 113
$ws.Range("H113").Value = 23857
$ws.Range("I113").Value = 20400
$ws.Range("J113").Value = 32499.5
$ws.Range("K113").Value = 20400
$ws.Range("L113").Value = 32499.5
$ws.Range("M113").Value = -17146
$ws.Range("N113").Value = -39007.5
# Row 136
$ws.Range("H136").Value = 67584.60000000001
$ws.Range("J136").Value = 67584.60000000001
$ws.Range("L136").Value = 67584.60000000001
$ws.Range("N136").Value = -77784.60000000001
# Row 138
$ws.Range("H138").Value = 2363.0505
$ws.Range("I138").Value = 1248.9286
$ws.Range("J138").Value = 2546.553
$ws.Range("K138").Value = 3746.7858
$ws.Range("L138").Value = 7639.659
$ws.Range("M138").Value = 1393.2142
$ws.Range("N138").Value = -17919.659

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8803.384
$ws.Range("I32").Value = 1938.2727
$ws.Range("K32").Value = 1938.2727
$ws.Range("M32").Value = -1651.2727
# Row 74
$ws.Range("H74").Value = 10331.48
$ws.Range("I74").Value = 2198.3333
$ws.Range("K74").Value = 2198.3333
$ws.Range("M74").Value = -1324.3333
# Row 77
$ws.Range("H77").Value = 10331.48
$ws.Range("I77").Value = 2198.3333
$ws.Range("K77").Value = 10991.6665
$ws.Range("M77").Value = -6623.666499999999
# Row 140
$ws.Range("H140").Value = 79663.336
$ws.Range("J140").Value = 79663.336
$ws.Range("L140").Value = 79663.336
$ws.Range("N140").Value = -90023.336

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 507.33334
$ws.Range("I22").Value = 509.5
$ws.Range("K22").Value = 509.5
$ws.Range("M22").Value = -336.5
# Row 24
$ws.Range("H24").Value = 2210.5
$ws.Range("I24").Value = 1152.6
$ws.Range("J24").Value = 7500
$ws.Range("K24").Value = 1152.6
$ws.Range("L24").Value = 7500
$ws.Range("M24").Value = -917.5999999999999
$ws.Range("N24").Value = -7970
# Row 86
$ws.Range("H86").Value = 1855.4166
$ws.Range("I86").Value = 1796.25
$ws.Range("J86").Value = 1973.75
$ws.Range("K86").Value = 1796.25
$ws.Range("L86").Value = 1973.75
$ws.Range("M86").Value = -673.25
$ws.Range("N86").Value = -4219.75
# Row 89
$ws.Range("H89").Value = 1855.4166
$ws.Range("I89").Value = 1796.25
$ws.Range("J89").Value = 1973.75
$ws.Range("K89").Value = 8981.25
$ws.Range("L89").Value = 9868.75
$ws.Range("M89").Value = -3365.25
$ws.Range("N89").Value = -21100.75
# Row 104
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
# Row 105
$ws.Range("H105").Value = 1273.7188
$ws.Range("I105").Value = 665.3158
$ws.Range("K105").Value = 665.3158
$ws.Range("M105").Value = 1081.6842

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9960.929
$ws.Range("I31").Value = 5899.357
$ws.Range("J31").Value = 14022.5
$ws.Range("K31").Value = 5899.357
$ws.Range("L31").Value = 14022.5
$ws.Range("M31").Value = -5604.357
$ws.Range("N31").Value = -14612.5
# Row 34
$ws.Range("H34").Value = 9960.929
$ws.Range("I34").Value = 5899.357
$ws.Range("J34").Value = 14022.5
$ws.Range("K34").Value = 5899.357
$ws.Range("L34").Value = 14022.5
$ws.Range("M34").Value = -5697.357
$ws.Range("N34").Value = -14426.5
# Row 58
$ws.Range("H58").Value = 11836.575
$ws.Range("I58").Value = 4966.1304
$ws.Range("J58").Value = 21131.883
$ws.Range("K58").Value = 4966.1304
$ws.Range("L58").Value = 21131.883
$ws.Range("M58").Value = -4763.1304
$ws.Range("N58").Value = -21537.883
# Row 99
$ws.Range("H99").Value = 11194.223
$ws.Range("I99").Value = 4916.3335
$ws.Range("K99").Value = 4916.3335
$ws.Range("M99").Value = -3418.3335
# Row 105
$ws.Range("H105").Value = 12485.23
$ws.Range("I105").Value = 13101.125
$ws.Range("J105").Value = 11499.8
$ws.Range("K105").Value = 13101.125
$ws.Range("L105").Value = 11499.8
$ws.Range("M105").Value = -11354.125
$ws.Range("N105").Value = -14993.8
# Row 126
$ws.Range("H126").Value = 11194.223
$ws.Range("I126").Value = 4916.3335
$ws.Range("K126").Value = 14749.0005
$ws.Range("M126").Value = -12279.0005
# Row 136
$ws.Range("H136").Value = 11836.575
$ws.Range("I136").Value = 4966.1304
$ws.Range("J136").Value = 21131.883
$ws.Range("K136").Value = 14898.3912
$ws.Range("L136").Value = 63395.649
$ws.Range("M136").Value = -12348.3912
$ws.Range("N136").Value = -68495.649

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 4920843.5
$ws.Range("I5").Value = 1826.6471
$ws.Range("J5").Value = 9839861
$ws.Range("K5").Value = 5479.9413
$ws.Range("L5").Value = 29519583
$ws.Range("M5").Value = -5367.9413
$ws.Range("N5").Value = -29519807
# Row 38
$ws.Range("H38").Value = 104.125
$ws.Range("I38").Value = 155.16667
$ws.Range("J38").Value = 73.5
$ws.Range("K38").Value = 465.50001
$ws.Range("L38").Value = 220.5
$ws.Range("M38").Value = -118.50001
$ws.Range("N38").Value = -914.5
# Row 56
$ws.Range("H56").Value = 166673140
$ws.Range("I56").Value = 166673140
$ws.Range("K56").Value = 166673140
$ws.Range("M56").Value = -166672610
# Row 113
$ws.Range("H113").Value = 1247.2941
$ws.Range("J113").Value = 1508.5834
$ws.Range("L113").Value = 4525.7502
$ws.Range("N113").Value = -8865.7502
# Row 114
$ws.Range("H114").Value = 3108
$ws.Range("I114").Value = 2366.6667
$ws.Range("J114").Value = 3293.3333
$ws.Range("K114").Value = 7100.000100000001
$ws.Range("L114").Value = 9879.999899999999
$ws.Range("M114").Value = -3846.000100000001
$ws.Range("N114").Value = -16387.9999
# Row 131
$ws.Range("H131").Value = 1441.13
$ws.Range("I131").Value = 675.5714
$ws.Range("J131").Value = 1498.7527
$ws.Range("K131").Value = 2026.7142
$ws.Range("L131").Value = 4496.2581
$ws.Range("M131").Value = 3013.2858
$ws.Range("N131").Value = -14576.2581
# Row 135
$ws.Range("H135").Value = 4920843.5
$ws.Range("I135").Value = 1826.6471
$ws.Range("J135").Value = 9839861
$ws.Range("K135").Value = 16439.8239
$ws.Range("L135").Value = 88558749
$ws.Range("M135").Value = -13904.8239
$ws.Range("N135").Value = -88563819
# Row 140
$ws.Range("H140").Value = 1524
$ws.Range("I140").Value = 1455
$ws.Range("K140").Value = 4365
$ws.Range("M140").Value = 815

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 13252.18
$ws.Range("I80").Value = 11265.8
$ws.Range("J80").Value = 15343.105
$ws.Range("K80").Value = 11265.8
$ws.Range("L80").Value = 15343.105
$ws.Range("M80").Value = -10267.8
$ws.Range("N80").Value = -17339.105
# Row 83
$ws.Range("H83").Value = 13252.18
$ws.Range("I83").Value = 11265.8
$ws.Range("J83").Value = 15343.105
$ws.Range("K83").Value = 56329
$ws.Range("L83").Value = 76715.52499999999
$ws.Range("M83").Value = -51337
$ws.Range("N83").Value = -86699.52499999999
# Row 122
$ws.Range("H122").Value = 3042.3704
$ws.Range("I122").Value = 3172
$ws.Range("J122").Value = 2005.3334
$ws.Range("K122").Value = 9516
$ws.Range("L122").Value = 6016.0002
$ws.Range("M122").Value = -7066
$ws.Range("N122").Value = -10916.0002

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 18250.25
$ws.Range("I7").Value = 33334.668
$ws.Range("J7").Value = 13222.111
$ws.Range("K7").Value = 33334.668
$ws.Range("L7").Value = 13222.111
$ws.Range("M7").Value = -33222.668
$ws.Range("N7").Value = -13446.111
# Row 40
$ws.Range("H40").Value = 5084.085
$ws.Range("I40").Value = 2592.697
$ws.Range("J40").Value = 10956.643
$ws.Range("K40").Value = 2592.697
$ws.Range("L40").Value = 10956.643
$ws.Range("M40").Value = -2456.697
$ws.Range("N40").Value = -11228.643
# Row 61
$ws.Range("H61").Value = 3727.2727
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2298
# Row 113
$ws.Range("H113").Value = 3727.2727
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
# Row 126
$ws.Range("H126").Value = 18250.25
$ws.Range("I126").Value = 33334.668
$ws.Range("J126").Value = 13222.111
$ws.Range("K126").Value = 100004.004
$ws.Range("L126").Value = 39666.333
$ws.Range("M126").Value = -97534.00399999999
$ws.Range("N126").Value = -44606.333
# Row 140
$ws.Range("H140").Value = 177746.38
$ws.Range("J140").Value = 177746.38
$ws.Range("L140").Value = 177746.38
$ws.Range("N140").Value = -188106.38

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# Row 126
$ws.Range("H126").Value = 15600
$ws.Range("I126").Value = 15750
$ws.Range("K126").Value = 47250
$ws.Range("M126").Value = -44780
# Row 132
$ws.Range("H132").Value = 4794.3506
$ws.Range("I132").Value = 2200.28
$ws.Range("K132").Value = 6600.84
$ws.Range("M132").Value = -4070.84

Write-Output "edit complete"